# Fruta / hortaliza, semanal
# Inserts two new weekly price rows (Kiwi - Provincia de Curico, fecha 44522)
# right after row 394, shifting the existing rows 395-466 down to 397-468.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before the current row 395. This pushes the
# existing data (old rows 395-466) down to rows 397-468, matching the
# target dimension A1:T468.
$ws.Rows("395:396").Insert()

# --- New row 395 ---
$ws.Range("A395").Value = 9
$ws.Range("B395").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C395").Value = "Metropolitana"
$ws.Range("D395").Value = 44522
$ws.Range("E395").Value = 13
$ws.Range("F395").Value = "Fruta"
$ws.Range("G395").Value = 100101
$ws.Range("H395").Value = "Berries"
$ws.Range("I395").Value = 100101007
$ws.Range("J395").Value = "Kiwi"
$ws.Range("K395").Value = "Hayward"
$ws.Range("L395").Value = "Primera"
$ws.Range("M395").Value = 300
$ws.Range("N395").Value = 9000
$ws.Range("O395").Value = 9000
$ws.Range("P395").Value = 9000
$ws.Range("Q395").Value = "$/bandeja 10 kilos"
$ws.Range("R395").Value = "Provincia de Curicó"
$ws.Range("S395").Value = 900
$ws.Range("T395").Value = 10

# --- New row 396 ---
$ws.Range("A396").Value = 9
$ws.Range("B396").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C396").Value = "Metropolitana"
$ws.Range("D396").Value = 44522
$ws.Range("E396").Value = 13
$ws.Range("F396").Value = "Fruta"
$ws.Range("G396").Value = 100101
$ws.Range("H396").Value = "Berries"
$ws.Range("I396").Value = 100101007
$ws.Range("J396").Value = "Kiwi"
$ws.Range("K396").Value = "Hayward"
$ws.Range("L396").Value = "Segunda"
$ws.Range("M396").Value = 250
$ws.Range("N396").Value = 7000
$ws.Range("O396").Value = 7000
$ws.Range("P396").Value = 7000
$ws.Range("Q396").Value = "$/bandeja 10 kilos"
$ws.Range("R396").Value = "Provincia de Curicó"
$ws.Range("S396").Value = 700
$ws.Range("T396").Value = 10

# Make sure the Date column keeps the date/time number format used
# throughout column D (style index 2 in the original workbook).
$ws.Range("D395").NumberFormat = $ws.Range("D394").NumberFormat
$ws.Range("D396").NumberFormat = $ws.Range("D394").NumberFormat
